$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = 2,3,5,6,8,10,11,12,13,14,15,17,18,19,20,21,22,24,26,28,29,31,32,34,36,37,38,39,40,41,43,44,45,46,47,48,50,52,54,55,57,58,60,62,63,64,65,66,67,69,70,71,72,73,74,76,78,80,81,82,83,84,85,86,87,90,92,93,94,96,99,101,106,107,108,109,110,111,112,113,116,118,119,120,122,125,127,132,133,134,135,136,137,138,139,142,144,145,146,148,151,153

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)
    $value = $cell.Value()
    if ($value) {
        $parts = $value -split ', '
        $reversed = @()
        for ($i = $parts.Count - 1; $i -ge 0; $i--) {
            $reversed += $parts[$i]
        }
        $cell.Value = [string]::Join(', ', $reversed)
    }
}
